# ajout d'une url pour la commission
# Insert two new columns (G:H) before the former "Fax" column, then
# populate the new header/data cells for the commission's Url and the
# "show mail" flag used by the mail display column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank columns at G:H - existing columns G.. shift right to I..
# Excel clones the formatting of the column to the left (F) for the new
# columns, which already yields the correct row-3 blank style (s=7).
$ws.Columns("G:H").Insert()

# Match the original authored column widths as closely as the ColumnWidth
# property (character-width units, 1/6 px granularity) allows.
$ws.Columns("G").ColumnWidth = 18.6666666666667
$ws.Columns("H").ColumnWidth = 23.8333333333333

# Populate header row (row 1) - order matches the shared-string creation
# order from the authored workbook (Url, then url formula, then
# show-mail formula, then the "Affichage mail" header).
$ws.Range("H1").Value = "Url"
$ws.Range("H2").Value = '${com.url}'
$ws.Range("G2").Value = '${com.temShowMail}'
$ws.Range("G1").Value = "Affichage mail"

# Restore the author's last on-sheet selection.
$ws.Range("H8").Select()
